# Slide 9, shape 2 ("Прямоугольник: скругленные углы 4") currently has a
# single Uzbek run. The edit rewrites it as three runs:
#   1) ru-RU  "Доктор рассмотрит анализ клиента и при необходимости даст
#              ему инструкции и рекомендации. "
#   2) uz-Cyrl-UZ "Здесь будет кнопка для добавления "
#   3) ru-RU  "рекомендации."
#
# The fake TextRange.LanguageID setter in this COM shim always stamps the
# *first* run of the shape's text, no matter which sub-range object it is
# invoked on. So runs are built back-to-front with InsertBefore, fixing
# each run's language immediately after it becomes run #1.
#
# TextRange.Delete() (rather than Text = "") is used to clear the original
# text because it also drops the leftover <a:endParaRPr>, matching the
# target XML which no longer has one.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$run1Text = "Доктор рассмотрит анализ клиента и при необходимости даст ему инструкции и рекомендации. "
$run2Text = "Здесь будет кнопка для добавления "
$run3Text = "рекомендации."

# Start from a clean paragraph (no stray endParaRPr).
$tr.Delete()

# Run 3 (last) goes in first, then runs 2 and 1 are prepended in turn so
# each one is, at the moment its language is set, the first (and only)
# run seen by LanguageID.
$tr.Text = $run3Text
$tr.LanguageID = "ru-RU"

$run2 = $tr.InsertBefore($run2Text)
$run2.LanguageID = "uz-Cyrl-UZ"

$run1 = $tr.InsertBefore($run1Text)
$run1.LanguageID = "ru-RU"
